$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'308.45"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'3.42%"
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 7).Value = "'2"
$ws.Cells.Item(2, 7).Style = "Normal"

$ws.Cells.Item(3, 4).Value = "'36.42"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'4.04%"
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 7).Value = "'2"
$ws.Cells.Item(3, 7).Style = "Normal"

$ws.Cells.Item(4, 4).Value = "'5.143"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'4.00%"
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(4, 7).Value = "'2"
$ws.Cells.Item(4, 7).Style = "Normal"

$ws.Cells.Item(5, 4).Value = "'0.08164"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'4.74%"
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 7).Value = "'2"
$ws.Cells.Item(5, 7).Style = "Normal"

$ws.Cells.Item(6, 4).Value = "'1.940"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'1.06%"
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(6, 7).Value = "'2"
$ws.Cells.Item(6, 7).Style = "Normal"

$ws.Cells.Item(7, 4).Value = "'7.786"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'1.05%"
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(7, 7).Value = "'2"
$ws.Cells.Item(7, 7).Style = "Normal"

$ws.Cells.Item(8, 4).Value = "'0.9324"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'1.46%"
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(8, 7).Value = "'2"
$ws.Cells.Item(8, 7).Style = "Normal"

$ws.Cells.Item(9, 4).Value = "'0.1395"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'30.40%"
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(9, 7).Value = "'2"
$ws.Cells.Item(9, 7).Style = "Normal"

$ws.Cells.Item(10, 4).Value = "'0.1938"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'7.67%"
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(10, 7).Value = "'2"
$ws.Cells.Item(10, 7).Style = "Normal"

$ws.Cells.Item(11, 4).Value = "'0.09258"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'0.55%"
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(11, 7).Value = "'2"
$ws.Cells.Item(11, 7).Style = "Normal"

$ws.Cells.Item(12, 4).Value = "'0.03570"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'0.53%"
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(12, 7).Value = "'2"
$ws.Cells.Item(12, 7).Style = "Normal"

$ws.Cells.Item(13, 4).Value = "'0.09864"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'-0.06%"
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(13, 7).Value = "'2"
$ws.Cells.Item(13, 7).Style = "Normal"

$ws.Cells.Item(14, 4).Value = "'0.001412"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'1.69%"
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(14, 7).Value = "'2"
$ws.Cells.Item(14, 7).Style = "Normal"

$ws.Cells.Item(15, 4).Value = "'0.005965"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'3.21%"
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(15, 7).Value = "'2"
$ws.Cells.Item(15, 7).Style = "Normal"

$ws.Cells.Item(16, 4).Value = "'3.540"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'2.08%"
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(16, 7).Value = "'2"
$ws.Cells.Item(16, 7).Style = "Normal"

$ws.Cells.Item(17, 4).Value = "'4.172"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'4.09%"
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(17, 7).Value = "'2"
$ws.Cells.Item(17, 7).Style = "Normal"

$ws.Cells.Item(18, 4).Value = "'2.986"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'2.64%"
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(18, 7).Value = "'2"
$ws.Cells.Item(18, 7).Style = "Normal"

$ws.Cells.Item(19, 4).Value = "'0.3440"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'0.04%"
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(19, 7).Value = "'2"
$ws.Cells.Item(19, 7).Style = "Normal"

$ws.Cells.Item(20, 4).Value = "'0.1336"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'3.50%"
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(20, 7).Value = "'2"
$ws.Cells.Item(20, 7).Style = "Normal"

$ws.Cells.Item(21, 4).Value = "'4.888"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'-3.03%"
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(21, 7).Value = "'2"
$ws.Cells.Item(21, 7).Style = "Normal"

$ws.Cells.Item(22, 4).Value = "'0.2403"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'9.75%"
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(22, 7).Value = "'2"
$ws.Cells.Item(22, 7).Style = "Normal"

$ws.Cells.Item(23, 4).Value = "'0.04492"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'-1.03%"
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(23, 7).Value = "'2"
$ws.Cells.Item(23, 7).Style = "Normal"

$ws.Cells.Item(24, 2).Value = "BitKan"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(24, 4).Value = "'0.001211"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'-0.10%"
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(24, 7).Value = "'2"
$ws.Cells.Item(24, 7).Style = "Normal"

$ws.Cells.Item(25, 2).Value = "HotbitToken"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(25, 4).Value = "'0.004910"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'6.93%"
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(25, 7).Value = "'2"
$ws.Cells.Item(25, 7).Style = "Normal"

$ws.Cells.Item(26, 5).Value = "'-0.96%"
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(26, 7).Value = "'2"
$ws.Cells.Item(26, 7).Style = "Normal"

$ws.Cells.Item(27, 4).Value = "'0.0004442"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'5.92%"
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(27, 7).Value = "'2"
$ws.Cells.Item(27, 7).Style = "Normal"

$ws.Cells.Item(28, 7).Value = "'2"
$ws.Cells.Item(28, 7).Style = "Normal"

$ws.Cells.Item(29, 7).Value = "'2"
$ws.Cells.Item(29, 7).Style = "Normal"

$ws.Cells.Item(30, 7).Value = "'2"
$ws.Cells.Item(30, 7).Style = "Normal"

$ws.Cells.Item(31, 7).Value = "'2"
$ws.Cells.Item(31, 7).Style = "Normal"

$ws.Cells.Item(32, 7).Value = "'2"
$ws.Cells.Item(32, 7).Style = "Normal"

$ws.Cells.Item(33, 7).Value = "'2"
$ws.Cells.Item(33, 7).Style = "Normal"

$ws.Cells.Item(34, 7).Value = "'2"
$ws.Cells.Item(34, 7).Style = "Normal"

$ws.Cells.Item(35, 7).Value = "'2"
$ws.Cells.Item(35, 7).Style = "Normal"

$ws.Cells.Item(36, 7).Value = "'2"
$ws.Cells.Item(36, 7).Style = "Normal"

$ws.Cells.Item(37, 7).Value = "'2"
$ws.Cells.Item(37, 7).Style = "Normal"

$ws.Cells.Item(38, 7).Value = "'2"
$ws.Cells.Item(38, 7).Style = "Normal"

$ws.Cells.Item(39, 4).Value = "'0.02003"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'7.38%"
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(39, 7).Value = "'2"
$ws.Cells.Item(39, 7).Style = "Normal"

$ws.Cells.Item(40, 4).Value = "'0.04934"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'6.16%"
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(40, 7).Value = "'2"
$ws.Cells.Item(40, 7).Style = "Normal"

$ws.Cells.Item(41, 4).Value = "'0.01116"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'16.16%"
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(41, 7).Value = "'2"
$ws.Cells.Item(41, 7).Style = "Normal"

$ws.Cells.Item(42, 4).Value = "'0.007641"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'0.93%"
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(42, 7).Value = "'2"
$ws.Cells.Item(42, 7).Style = "Normal"

$ws.Cells.Item(43, 4).Value = "'0.1384"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'5.02%"
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(43, 7).Value = "'2"
$ws.Cells.Item(43, 7).Style = "Normal"

$ws.Cells.Item(44, 4).Value = "'0.002098"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'-0.84%"
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(44, 7).Value = "'2"
$ws.Cells.Item(44, 7).Style = "Normal"

$ws.Cells.Item(45, 4).Value = "'0.01061"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'-4.04%"
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(45, 7).Value = "'2"
$ws.Cells.Item(45, 7).Style = "Normal"

$ws.Cells.Item(46, 4).Value = "'0.00006446"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'7.04%"
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(46, 7).Value = "'2"
$ws.Cells.Item(46, 7).Style = "Normal"

$ws.Cells.Item(47, 4).Value = "'0.00000000749"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'-0.15%"
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(47, 7).Value = "'2"
$ws.Cells.Item(47, 7).Style = "Normal"

$ws.Cells.Item(48, 7).Value = "'2"
$ws.Cells.Item(48, 7).Style = "Normal"

$ws.Cells.Item(49, 4).Value = "'0.001191"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'-8.70%"
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(49, 7).Value = "'2"
$ws.Cells.Item(49, 7).Style = "Normal"

$ws.Cells.Item(50, 4).Value = "'0.00002098"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'-0.15%"
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(50, 7).Value = "'2"
$ws.Cells.Item(50, 7).Style = "Normal"

$ws.Cells.Item(51, 4).Value = "'0.0001998"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'-0.15%"
$ws.Cells.Item(51, 5).Style = "Normal"
$ws.Cells.Item(51, 7).Value = "'2"
$ws.Cells.Item(51, 7).Style = "Normal"
